$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C14: it was stored as text "8894" but should be a real number.
$ws.Range("C14").Value = 8894

# Add new row 15 with data for a new company record.
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "razon social 1"
# C15's "12545" must stay text (unlike C14 above). A leading apostrophe
# forces Excel to keep a numeric-looking literal as text instead of
# auto-coercing it to a number; reset the style afterwards so it doesn't
# leave a lingering quote-prefix number format on the cell.
$ws.Range("C15").Value = "'12545"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = 4330
$ws.Range("E15").Value = "Terminación y acabado de edificios"
$ws.Range("F15").Value = "representante 1"
$ws.Range("G15").Value = "correo 1"
$ws.Range("H15").Value = "telefono 1"
